# Refresh market/profit data (Universalis price pull) across the Leve profit sheets.
# Columns H-N hold currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# figures that are recomputed by the scheduled runner and written back as static values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 519.5294
$ws.Range("I33").Value = 592.2857
$ws.Range("J33").Value = 180
$ws.Range("K33").Value = 592.2857
$ws.Range("L33").Value = 180
$ws.Range("M33").Value = -363.2857
$ws.Range("N33").Value = -638
$ws.Range("H112").Value = 14595.263
$ws.Range("J112").Value = 14595.263
$ws.Range("L112").Value = 43785.789
$ws.Range("N112").Value = -46001.789
$ws.Range("H129").Value = 1328
$ws.Range("J129").Value = 1688.8846
$ws.Range("L129").Value = 5066.6538
$ws.Range("N129").Value = -15066.6538
$ws.Range("H141").Value = 5578.385
$ws.Range("I141").Value = 2714.889
$ws.Range("K141").Value = 8144.667
$ws.Range("M141").Value = -2964.667
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1966.6666
$ws.Range("I8").Value = 300
$ws.Range("J8").Value = 2800
$ws.Range("K8").Value = 300
$ws.Range("L8").Value = 2800
$ws.Range("M8").Value = -156
$ws.Range("N8").Value = -3088
$ws.Range("H10").Value = 33285
$ws.Range("J10").Value = 47430
$ws.Range("L10").Value = 47430
$ws.Range("N10").Value = -47770
$ws.Range("H11").Value = 1600300
$ws.Range("I11").Value = 2000225
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 2000225
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -2000081
$ws.Range("N11").Value = -888
$ws.Range("H12").Value = 1068
$ws.Range("J12").Value = 1302
$ws.Range("L12").Value = 1302
$ws.Range("N12").Value = -1648
$ws.Range("H13").Value = 895
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 895
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 895
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1183
$ws.Range("H14").Value = 1325
$ws.Range("J14").Value = 1325
$ws.Range("L14").Value = 1325
$ws.Range("N14").Value = -1675
$ws.Range("H63").Value = 3726.25
$ws.Range("I63").Value = 2010.9
$ws.Range("K63").Value = 2010.9
$ws.Range("M63").Value = -1324.9
$ws.Range("H66").Value = 3726.25
$ws.Range("I66").Value = 2010.9
$ws.Range("K66").Value = 10054.5
$ws.Range("M66").Value = -6622.5
$ws.Range("H74").Value = 10418694
$ws.Range("I74").Value = 1244.0286
$ws.Range("J74").Value = 38465670
$ws.Range("K74").Value = 1244.0286
$ws.Range("L74").Value = 38465670
$ws.Range("M74").Value = -370.0286000000001
$ws.Range("N74").Value = -38467418
$ws.Range("H77").Value = 10418694
$ws.Range("I77").Value = 1244.0286
$ws.Range("J77").Value = 38465670
$ws.Range("K77").Value = 6220.143
$ws.Range("L77").Value = 192328350
$ws.Range("M77").Value = -1852.143
$ws.Range("N77").Value = -192337086
$ws.Range("H80").Value = 17985.375
$ws.Range("J80").Value = 17985.375
$ws.Range("L80").Value = 17985.375
$ws.Range("N80").Value = -19981.375
$ws.Range("H83").Value = 17985.375
$ws.Range("J83").Value = 17985.375
$ws.Range("L83").Value = 53956.125
$ws.Range("N83").Value = -63940.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13181.546
$ws.Range("I82").Value = 4153.5
$ws.Range("J82").Value = 18340.428
$ws.Range("K82").Value = 4153.5
$ws.Range("L82").Value = 18340.428
$ws.Range("M82").Value = -3770.5
$ws.Range("N82").Value = -19106.428
$ws.Range("H85").Value = 13181.546
$ws.Range("I85").Value = 4153.5
$ws.Range("J85").Value = 18340.428
$ws.Range("K85").Value = 4153.5
$ws.Range("L85").Value = 18340.428
$ws.Range("M85").Value = -2827.5
$ws.Range("N85").Value = -20992.428
$ws.Range("H107").Value = 335300
$ws.Range("I107").Value = 501450
$ws.Range("K107").Value = 501450
$ws.Range("M107").Value = -499530
$ws.Range("H134").Value = 2100.5186
$ws.Range("I134").Value = 1828.5714
$ws.Range("J134").Value = 3052.3333
$ws.Range("K134").Value = 5485.7142
$ws.Range("L134").Value = 9156.999899999999
$ws.Range("M134").Value = -2950.7142
$ws.Range("N134").Value = -14226.9999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1456.5
$ws.Range("I5").Value = 114.84615
$ws.Range("J5").Value = 3394.4443
$ws.Range("K5").Value = 114.84615
$ws.Range("L5").Value = 3394.4443
$ws.Range("M5").Value = -2.846149999999994
$ws.Range("N5").Value = -3618.4443
$ws.Range("H31").Value = 5307.931
$ws.Range("I31").Value = 1427.8235
$ws.Range("J31").Value = 6916.756
$ws.Range("K31").Value = 1427.8235
$ws.Range("L31").Value = 6916.756
$ws.Range("M31").Value = -1132.8235
$ws.Range("N31").Value = -7506.756
$ws.Range("H34").Value = 5307.931
$ws.Range("I34").Value = 1427.8235
$ws.Range("J34").Value = 6916.756
$ws.Range("K34").Value = 1427.8235
$ws.Range("L34").Value = 6916.756
$ws.Range("M34").Value = -1225.8235
$ws.Range("N34").Value = -7320.756
$ws.Range("H50").Value = 12999.167
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 12999.167
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 12999.167
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -14249.167
$ws.Range("H51").Value = 17999.2
$ws.Range("J51").Value = 17999.2
$ws.Range("L51").Value = 17999.2
$ws.Range("N51").Value = -19471.2
$ws.Range("H57").Value = 39400
$ws.Range("J57").Value = 39400
$ws.Range("L57").Value = 39400
$ws.Range("N57").Value = -40520
$ws.Range("H59").Value = 14911.5
$ws.Range("J59").Value = 14911.5
$ws.Range("L59").Value = 14911.5
$ws.Range("N59").Value = -17201.5
$ws.Range("H60").Value = 16309.818
$ws.Range("I60").Value = 8000
$ws.Range("J60").Value = 17140.8
$ws.Range("K60").Value = 8000
$ws.Range("L60").Value = 17140.8
$ws.Range("M60").Value = -7489
$ws.Range("N60").Value = -18162.8
$ws.Range("H61").Value = 17999.2
$ws.Range("J61").Value = 17999.2
$ws.Range("L61").Value = 17999.2
$ws.Range("N61").Value = -18695.2
$ws.Range("H62").Value = 4190
$ws.Range("I62").Value = 3987.5
$ws.Range("K62").Value = 3987.5
$ws.Range("M62").Value = -3363.5
$ws.Range("H65").Value = 4190
$ws.Range("I65").Value = 3987.5
$ws.Range("K65").Value = 19937.5
$ws.Range("M65").Value = -16817.5
$ws.Range("H68").Value = 22381.75
$ws.Range("J68").Value = 22381.75
$ws.Range("L68").Value = 22381.75
$ws.Range("N68").Value = -23879.75
$ws.Range("H71").Value = 22381.75
$ws.Range("J71").Value = 22381.75
$ws.Range("L71").Value = 67145.25
$ws.Range("N71").Value = -74633.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 231041.77
$ws.Range("I11").Value = 1500095.5
$ws.Range("J11").Value = 304.72726
$ws.Range("K11").Value = 4500286.5
$ws.Range("L11").Value = 914.18178
$ws.Range("M11").Value = -4500146.5
$ws.Range("N11").Value = -1194.18178
$ws.Range("H132").Value = 2147.8484
$ws.Range("I132").Value = 2164.4211
$ws.Range("J132").Value = 2125.3572
$ws.Range("K132").Value = 19479.7899
$ws.Range("L132").Value = 19128.2148
$ws.Range("M132").Value = -16949.7899
$ws.Range("N132").Value = -24188.2148
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1240
$ws.Range("I102").Value = 1250
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 1250
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 372
$ws.Range("N102").Value = -4444
$ws.Range("H113").Value = 80530.92999999999
$ws.Range("I113").Value = 86494.84
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 86494.84
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -84324.84
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 11527.909
$ws.Range("I122").Value = 16743.857
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 50231.571
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -47781.571
$ws.Range("N122").Value = -12100
$ws.Range("H123").Value = 8322
$ws.Range("J123").Value = 8322
$ws.Range("L123").Value = 8322
$ws.Range("N123").Value = -13222
$ws.Range("H132").Value = 2938.6924
$ws.Range("I132").Value = 2133.625
$ws.Range("J132").Value = 4226.8
$ws.Range("K132").Value = 6400.875
$ws.Range("L132").Value = 12680.4
$ws.Range("M132").Value = -3870.875
$ws.Range("N132").Value = -17740.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 81753.5
$ws.Range("J14").Value = 81753.5
$ws.Range("L14").Value = 81753.5
$ws.Range("N14").Value = -82097.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3500
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3772
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 40000
$ws.Range("J27").Value = 40000
$ws.Range("L27").Value = 40000
$ws.Range("N27").Value = -40138
$ws.Range("H122").Value = 1821.1666
$ws.Range("I122").Value = 1878.2858
$ws.Range("J122").Value = 1621.25
$ws.Range("K122").Value = 5634.857400000001
$ws.Range("L122").Value = 4863.75
$ws.Range("M122").Value = -3184.857400000001
$ws.Range("N122").Value = -9763.75
$ws.Range("H132").Value = 6208215
$ws.Range("I132").Value = 2338.75
$ws.Range("J132").Value = 19447418
$ws.Range("K132").Value = 7016.25
$ws.Range("L132").Value = 58342254
$ws.Range("M132").Value = -4486.25
$ws.Range("N132").Value = -58347314
